# Refined metadata to be additional tab
#
# 1. Update the F-column ("time_taken") timestamps on the "data" sheet to
#    reflect the re-run query time.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself, mirroring the header style used on "data".

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken values on the data sheet -------------------
$timestamps = @(
    "2021-10-05 14:20:21.204791",
    "2021-10-05 14:20:21.204799",
    "2021-10-05 14:20:21.204802",
    "2021-10-05 14:20:21.204804",
    "2021-10-05 14:20:21.204807",
    "2021-10-05 14:20:21.204810",
    "2021-10-05 14:20:21.204812",
    "2021-10-05 14:20:21.204814",
    "2021-10-05 14:20:21.204817",
    "2021-10-05 14:20:21.204820",
    "2021-10-05 14:20:21.204822",
    "2021-10-05 14:20:21.204825",
    "2021-10-05 14:20:21.204827",
    "2021-10-05 14:20:21.204830",
    "2021-10-05 14:20:21.204832",
    "2021-10-05 14:20:21.204835",
    "2021-10-05 14:20:21.204837",
    "2021-10-05 14:20:21.204840",
    "2021-10-05 14:20:21.204842",
    "2021-10-05 14:20:21.204845",
    "2021-10-05 14:20:21.204847",
    "2021-10-05 14:20:21.204850",
    "2021-10-05 14:20:21.204852",
    "2021-10-05 14:20:21.204855",
    "2021-10-05 14:20:21.204858",
    "2021-10-05 14:20:21.204860",
    "2021-10-05 14:20:21.204863",
    "2021-10-05 14:20:21.204865",
    "2021-10-05 14:20:21.204868",
    "2021-10-05 14:20:21.204870",
    "2021-10-05 14:20:21.204872",
    "2021-10-05 14:20:21.204875",
    "2021-10-05 14:20:21.204878",
    "2021-10-05 14:20:21.204880",
    "2021-10-05 14:20:21.204883",
    "2021-10-05 14:20:21.204885",
    "2021-10-05 14:20:21.204887",
    "2021-10-05 14:20:21.204890",
    "2021-10-05 14:20:21.204892",
    "2021-10-05 14:20:21.204895",
    "2021-10-05 14:20:21.204897",
    "2021-10-05 14:20:21.204900",
    "2021-10-05 14:20:21.204902",
    "2021-10-05 14:20:21.204905",
    "2021-10-05 14:20:21.204907",
    "2021-10-05 14:20:21.204910",
    "2021-10-05 14:20:21.204912",
    "2021-10-05 14:20:21.204915",
    "2021-10-05 14:20:21.204917",
    "2021-10-05 14:20:21.204919",
    "2021-10-05 14:20:21.204922",
    "2021-10-05 14:20:21.204924",
    "2021-10-05 14:20:21.204927",
    "2021-10-05 14:20:21.204930",
    "2021-10-05 14:20:21.204932",
    "2021-10-05 14:20:21.204934",
    "2021-10-05 14:20:21.204937",
    "2021-10-05 14:20:21.204939",
    "2021-10-05 14:20:21.204942",
    "2021-10-05 14:20:21.204944",
    "2021-10-05 14:20:21.204947",
    "2021-10-05 14:20:21.204949",
    "2021-10-05 14:20:21.204951",
    "2021-10-05 14:20:21.204954",
    "2021-10-05 14:20:21.204957",
    "2021-10-05 14:20:21.204960",
    "2021-10-05 14:20:21.204962",
    "2021-10-05 14:20:21.204965",
    "2021-10-05 14:20:21.204967",
    "2021-10-05 14:20:21.204970",
    "2021-10-05 14:20:21.204972",
    "2021-10-05 14:20:21.204974",
    "2021-10-05 14:20:21.204977",
    "2021-10-05 14:20:21.204979"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. Add the metadata worksheet ----------------------------------------
$metadata = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$metadata.Name = "metadata"

# Header row + the A2 index cell reuse the same bold/bordered/centered
# header style already used on "data" (style index 1) -- copy the format
# from an existing header cell instead of building a new style combo.
$data.Range("B1").Copy()
$metadata.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$metadata.Range("B1").Value = "data_name"
$metadata.Range("C1").Value = "data_id"
$metadata.Range("D1").Value = "data_version"
$metadata.Range("E1").Value = "data_version_created"
$metadata.Range("F1").Value = "panel_query_time"
$metadata.Range("G1").Value = "panel_get_request"

# Data row
$metadata.Range("A2").Value = 0
$metadata.Range("B2").Value = "Familial pulmonary fibrosis"
$metadata.Range("C2").Value = 200

# "1.16" is stored as text (matching the source export), not the float
# 1.16 -- force text entry, then strip the number-format override back off
# so the cell ends up on the plain default style, just like a cell that was
# never explicitly formatted.
$metadata.Range("D2").NumberFormat = "@"
$metadata.Range("D2").Value = "1.16"
$metadata.Range("I1").Copy()
$metadata.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$metadata.Range("E2").Value = "2021-05-10T15:43:20.944111Z"
$metadata.Range("F2").Value = "2021-10-05 14:20:21.201131"
$metadata.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/200/?format=json"

$data.Select()
